$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.573030829429626
$ws.Range("B1").Value = 1.780390381813049
$ws.Range("C1").Value = 1.837659358978271
$ws.Range("D1").Value = 2.29645299911499
$ws.Range("E1").Value = 3.461513757705688
